# UPDATE technology portfolios for Norway
$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 8.5787499999999994
$ws2025.Range("C2").Value = 878666
$ws2025.Range("D2").Value = 360000

$sheetNames = @("2030", "2035", "2040", "2045", "2050")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B2").Value = 8.5787499999999994
}

$wb.Save()
